$p = $ppt.ActivePresentation
$s = $p.Slides.Item(31)
$sh = $s.Shapes.Item(2)

# Resize / reposition the content placeholder
$sh.Left = 53.33338737487793
$sh.Top = 112.2936220472441
$sh.Width = 676.9029921259844
$sh.Height = 395.00914001464844

$tr = $sh.TextFrame.TextRange
$tr.Text = "Redis is one of the most popular in-memory key-data-structure database.`rUse cases:`rCache`rMessaging - publisher/subscriber`rDistributed session/access token`rTemporary data with TTL - i.e. reservation, user block, voting system one per day`rReal-time access data - inventory, product prices.`rData structures:`rString - store values as string, ops: SET,GET, SETNX, INCR, etc`rList - list of strings, Ops: LPUSH, LPOP, LLEN, etc`rSet - unordered collection of strings with no repetition. SADD, SUNION, SINTER, etc`rSortedSet - ordered collection of strings with no repetition. ZADD, ZRANGEBYSCORE, etc`rHashe - is map between key and string value. HSET, HGET, HINCRBY, etc `rLimitations: `rValues up to 512MB`rCollections up to 2^32-1 ( more than 4 billion ) elements`r"

# Enable shrink-text-on-overflow autofit
$sh.TextFrame.AutoSize = 2

# Paragraph indent levels
$tr.Paragraphs(1,1).IndentLevel = 1
$tr.Paragraphs(2,1).IndentLevel = 1
$tr.Paragraphs(3,1).IndentLevel = 2
$tr.Paragraphs(4,1).IndentLevel = 2
$tr.Paragraphs(5,1).IndentLevel = 2
$tr.Paragraphs(6,1).IndentLevel = 2
$tr.Paragraphs(7,1).IndentLevel = 2
$tr.Paragraphs(8,1).IndentLevel = 1
$tr.Paragraphs(9,1).IndentLevel = 2
$tr.Paragraphs(10,1).IndentLevel = 2
$tr.Paragraphs(11,1).IndentLevel = 2
$tr.Paragraphs(12,1).IndentLevel = 2
$tr.Paragraphs(13,1).IndentLevel = 2
$tr.Paragraphs(14,1).IndentLevel = 1
$tr.Paragraphs(15,1).IndentLevel = 2
$tr.Paragraphs(16,1).IndentLevel = 2
$tr.Paragraphs(17,1).IndentLevel = 2

# Run-level formatting: underline + hyperlinks on "etc" references
$run = $tr.Characters(366,3)
$run.Font.Underline = 1
$run.ActionSettings.Item(1).Hyperlink.Address = "https://redis.io/commands/?group=string"
$run = $tr.Characters(418,3)
$run.Font.Underline = 1
$run.ActionSettings.Item(1).Hyperlink.Address = "https://redis.io/commands/?group=list"
$run = $tr.Characters(502,3)
$run.Font.Underline = 1
$run.ActionSettings.Item(1).Hyperlink.Address = "https://redis.io/commands/?group=set"
$run = $tr.Characters(589,3)
$run.Font.Underline = 1
$run.ActionSettings.Item(1).Hyperlink.Address = "https://redis.io/commands/?group=sorted_set"
$run = $tr.Characters(659,3)
$run.Font.Underline = 1
$run.ActionSettings.Item(1).Hyperlink.Address = "https://redis.io/commands/?group=hash"
